$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.906.25"
$ws.Range("E2").Value = "  +3.94%  "

$ws.Range("D3").Value = "2.489.60"
$ws.Range("E3").Value = "  +2.35%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.07"
$ws.Range("E5").Value = "  +1.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.85"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.520"
$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  +1.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.74"
$ws.Range("E10").Value = "  +2.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("E11").Value = "  +0.80%  "

$ws.Range("E12").Value = "  +0.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.25"
$ws.Range("E13").Value = "  -1.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.18"
$ws.Range("E14").Value = "  +3.07%  "

$ws.Range("D15").Value = "2.892.42"
$ws.Range("E15").Value = "  +2.99%  "

$ws.Range("D16").Value = "2.515.50"
$ws.Range("E16").Value = "  +3.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.840"
$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("D18").Value = "46.926.50"
$ws.Range("E18").Value = "  +4.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.60"
$ws.Range("E19").Value = "  +1.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.57"
$ws.Range("E20").Value = "  +3.09%  "

$ws.Range("D21").Value = "0.0₃0933"
$ws.Range("E21").Value = "  +1.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.52"
$ws.Range("E22").Value = "  +2.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.29"
$ws.Range("E23").Value = "  +2.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  +2.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  +1.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.11"
$ws.Range("E26").Value = "  +2.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("E28").Value = "  +4.12%  "

$ws.Range("E29").Value = "  +0.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.87"
$ws.Range("E30").Value = "  +3.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.133"
$ws.Range("E31").Value = "  +4.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.39"
$ws.Range("E32").Value = "  +0.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.56"
$ws.Range("E33").Value = "  -1.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.29"
$ws.Range("E34").Value = "  +1.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0775"
$ws.Range("E35").Value = "  +1.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  +0.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.93"
$ws.Range("E37").Value = "  +0.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.57"
$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.94"
$ws.Range("E39").Value = "  +2.11%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "122.99"
$ws.Range("E40").Value = "  -3.17%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.111"
$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +2.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.19"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0295"
$ws.Range("E44").Value = "  +1.49%  "

$ws.Range("D45").Value = "1.964.04"
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.98"
$ws.Range("E46").Value = "  +0.59%  "

$ws.Range("E47").Value = "  +0.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.78"
$ws.Range("E48").Value = "  -0.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.07"
$ws.Range("E49").Value = "  -1.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.39"
$ws.Range("E50").Value = "  +15.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.78"
$ws.Range("E51").Value = "  +3.67%  "
